$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (row 1, A1)
$ws.Range("A1").Value = "Datos actualizados a 15 de Agosto de 2020 a las 18:20"

# Update country data rows (new stats snapshot; a few rows also swap rank/name)
# Row 4: Estados Unidos -> Estados Unidos
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 5490948
$ws.Range("C4").Value = 14682
$ws.Range("D4").Value = 2877010
$ws.Range("E4").Value = 2442019
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 384
$ws.Range("H4").Value = 171919

# Row 5: Brasil -> Brasil
$ws.Range("A5").Value = "Brasil"
$ws.Range("B5").Value = 3282101
$ws.Range("C5").Value = 3206
$ws.Range("D5").Value = 2384302
$ws.Range("E5").Value = 791191
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 37
$ws.Range("H5").Value = 106608

# Row 6: India -> India
$ws.Range("A6").Value = "India"
$ws.Range("B6").Value = 2565755
$ws.Range("C6").Value = 40533
$ws.Range("D6").Value = 1843756
$ws.Range("E6").Value = 672358
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 507
$ws.Range("H6").Value = 49641

# Row 12: Chile -> Chile
$ws.Range("A12").Value = "Chile"
$ws.Range("B12").Value = 383902
$ws.Range("C12").Value = 1791
$ws.Range("D12").Value = 356951
$ws.Range("E12").Value = 16556
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 55
$ws.Range("H12").Value = 10395

# Row 15: Reino Unido -> Reino Unido
$ws.Range("A15").Value = "Reino Unido"
$ws.Range("B15").Value = 317379
$ws.Range("C15").Value = 1012
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = 41361

# Row 20: Italia -> Italia
$ws.Range("A20").Value = "Italia"
$ws.Range("B20").Value = 253438
$ws.Range("C20").Value = 629
$ws.Range("D20").Value = 203640
$ws.Range("E20").Value = 14406
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 4
$ws.Range("H20").Value = 35392

# Row 27: Canada -> Canada
$ws.Range("A27").Value = "Canada"
$ws.Range("B27").Value = 121760
$ws.Range("C27").Value = 108
$ws.Range("D27").Value = 108044
$ws.Range("E27").Value = 4695
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 9021

# Row 30: Ecuador -> Ecuador
$ws.Range("A30").Value = "Ecuador"
$ws.Range("B30").Value = 100688
$ws.Range("C30").Value = 1279
$ws.Range("D30").Value = 79354
$ws.Range("E30").Value = 15269
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 35
$ws.Range("H30").Value = 6065

# Row 35: China -> Republica Dominicana
$ws.Range("A35").Value = "Republica Dominicana"
$ws.Range("B35").Value = 85545
$ws.Range("C35").Value = 1057
$ws.Range("D35").Value = 51356
$ws.Range("E35").Value = 32751
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = 29
$ws.Range("H35").Value = 1438

# Row 36: Republica Dominicana -> China
$ws.Range("A36").Value = "China"
$ws.Range("B36").Value = 84808
$ws.Range("C36").Value = 22
$ws.Range("D36").Value = 79519
$ws.Range("E36").Value = 655
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 0
$ws.Range("H36").Value = 4634

# Row 74: Chequia -> Chequia
$ws.Range("A74").Value = "Chequia"
$ws.Range("B74").Value = 19818
$ws.Range("C74").Value = 125
$ws.Range("D74").Value = 13751
$ws.Range("E74").Value = 5672
$ws.Range("F74").Value = 0
$ws.Range("G74").Value = 1
$ws.Range("H74").Value = 395

# Row 83: Republica de Macedonia -> Republica de Macedonia
$ws.Range("A83").Value = "Republica de Macedonia"
$ws.Range("B83").Value = 12653
$ws.Range("C83").Value = 138
$ws.Range("D83").Value = 9123
$ws.Range("E83").Value = 2991
$ws.Range("F83").Value = 0
$ws.Range("G83").Value = 4
$ws.Range("H83").Value = 539

# Row 100: Albania -> Albania
$ws.Range("A100").Value = "Albania"
$ws.Range("B100").Value = 7260
$ws.Range("C100").Value = 143
$ws.Range("D100").Value = 3746
$ws.Range("E100").Value = 3289
$ws.Range("F100").Value = 0
$ws.Range("G100").Value = 6
$ws.Range("H100").Value = 225

# Row 101: Mauritania -> Grecia
$ws.Range("A101").Value = "Grecia"
$ws.Range("B101").Value = 6858
$ws.Range("C101").Value = 226
$ws.Range("D101").Value = 3804
$ws.Range("E101").Value = 2828
$ws.Range("F101").Value = 0
$ws.Range("G101").Value = 3
$ws.Range("H101").Value = 226

# Row 102: Grecia -> Mauritania
$ws.Range("A102").Value = "Mauritania"
$ws.Range("B102").Value = 6676
$ws.Range("C102").Value = 0
$ws.Range("D102").Value = 5889
$ws.Range("E102").Value = 630
$ws.Range("F102").Value = 0
$ws.Range("G102").Value = 0
$ws.Range("H102").Value = 157

# Row 115: Namibia -> Namibia
$ws.Range("A115").Value = "Namibia"
$ws.Range("B115").Value = 3726
$ws.Range("C115").Value = 0
$ws.Range("D115").Value = 2342
$ws.Range("E115").Value = 1353
$ws.Range("F115").Value = 0
$ws.Range("G115").Value = 0
$ws.Range("H115").Value = 31

# Row 118: Somalia -> Cuba
$ws.Range("A118").Value = "Cuba"
$ws.Range("B118").Value = 3292
$ws.Range("C118").Value = 63
$ws.Range("D118").Value = 2568
$ws.Range("E118").Value = 636
$ws.Range("F118").Value = 0
$ws.Range("G118").Value = 0
$ws.Range("H118").Value = 88

# Row 119: Cuba -> Somalia
$ws.Range("A119").Value = "Somalia"
$ws.Range("B119").Value = 3250
$ws.Range("C119").Value = 0
$ws.Range("D119").Value = 2268
$ws.Range("E119").Value = 889
$ws.Range("F119").Value = 0
$ws.Range("G119").Value = 0
$ws.Range("H119").Value = 93

# Row 145: Jordania -> Jordania
$ws.Range("A145").Value = "Jordania"
$ws.Range("B145").Value = 1339
$ws.Range("C145").Value = 10
$ws.Range("D145").Value = 1229
$ws.Range("E145").Value = 99
$ws.Range("F145").Value = 0
$ws.Range("G145").Value = 0
$ws.Range("H145").Value = 11

# Row 152: Niger -> Niger
$ws.Range("A152").Value = "Niger"
$ws.Range("B152").Value = 1165
$ws.Range("C152").Value = 4
$ws.Range("D152").Value = 1077
$ws.Range("E152").Value = 19
$ws.Range("F152").Value = 0
$ws.Range("G152").Value = 0
$ws.Range("H152").Value = 69

# Row 160: Santo Tome y Principe -> Lesoto
$ws.Range("A160").Value = "Lesoto"
$ws.Range("B160").Value = 903
$ws.Range("C160").Value = 19
$ws.Range("D160").Value = 271
$ws.Range("E160").Value = 607
$ws.Range("F160").Value = 0
$ws.Range("G160").Value = 0
$ws.Range("H160").Value = 25

# Row 161: Lesoto -> Santo Tome y Principe
$ws.Range("A161").Value = "Santo Tome y Principe"
$ws.Range("B161").Value = 885
$ws.Range("C161").Value = 0
$ws.Range("D161").Value = 809
$ws.Range("E161").Value = 61
$ws.Range("F161").Value = 0
$ws.Range("G161").Value = 0
$ws.Range("H161").Value = 15

# Row 213: Islas Malvinas -> Montserrat
$ws.Range("A213").Value = "Montserrat"
$ws.Range("B213").Value = 13
$ws.Range("C213").Value = 0
$ws.Range("D213").Value = 12
$ws.Range("E213").Value = 0
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 0
$ws.Range("H213").Value = 1

# Row 214: Montserrat -> Islas Malvinas
$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("B214").Value = 13
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 13
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 0
